$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 507
$ws.Range("F4").Value = 487
$ws.Range("F5").Value = 904
$ws.Range("F6").Value = 135
$ws.Range("F7").Value = 883
$ws.Range("F8").Value = 690
$ws.Range("F9").Value = 145
$ws.Range("F10").Value = 45
$ws.Range("F11").Value = 70
$ws.Range("F12").Value = 746
$ws.Range("F13").Value = 234
$ws.Range("F14").Value = 530
$ws.Range("F15").Value = 477
$ws.Range("F16").Value = 1260
$ws.Range("F17").Value = 108
$ws.Range("F18").Value = 56
$ws.Range("F19").Value = 999
$ws.Range("F20").Value = 2743
$ws.Range("G20").Value = 60
$ws.Range("F21").Value = 1216
$ws.Range("F22").Value = 625
$ws.Range("F23").Value = 154
$ws.Range("F24").Value = 1221
$ws.Range("F26").Value = 938
$ws.Range("F27").Value = 104
$ws.Range("F28").Value = 1266

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 501
$ws.Range("F4").Value = 343
$ws.Range("F8").Value = 35

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 711

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 711
$ws.Range("F3").Value = 507
$ws.Range("F6").Value = 487
$ws.Range("F7").Value = 501
$ws.Range("F8").Value = 501
$ws.Range("F9").Value = 343
$ws.Range("F12").Value = 904
$ws.Range("F13").Value = 135
$ws.Range("F14").Value = 883
$ws.Range("F15").Value = 690
$ws.Range("F16").Value = 145
$ws.Range("F18").Value = 45
$ws.Range("F19").Value = 35
$ws.Range("F22").Value = 70
$ws.Range("F24").Value = 746
$ws.Range("F25").Value = 234
$ws.Range("F26").Value = 530
$ws.Range("F27").Value = 477
$ws.Range("F28").Value = 1260
$ws.Range("F29").Value = 108
$ws.Range("F30").Value = 56
$ws.Range("F31").Value = 999
$ws.Range("F32").Value = 2743
$ws.Range("G32").Value = 60
$ws.Range("F33").Value = 1216
$ws.Range("F34").Value = 625
$ws.Range("F35").Value = 154
$ws.Range("F36").Value = 1221
$ws.Range("F39").Value = 938
$ws.Range("F40").Value = 104
$ws.Range("F41").Value = 1266
